# Menu-Languages.docx (Chinese Simplified) -- new translations
# This script reproduces three paragraph insertions / edits:
#   1. Insert a new "The Other Side" paragraph right after the first
#      "SmartCard" paragraph (before "Resources").
#   2. Insert two new paragraphs ("Exchanges Listing Guide" and a
#      duplicate "Graphics") right before the existing "Graphics"
#      paragraph, and turn that original "Graphics" paragraph's run
#      into "Marketing Materials" (Times New Roman).
#   3. Insert a new "Guide" paragraph right after "Electrum Wallet".

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------
# Change 1: "The Other Side" after the first "SmartCard" paragraph
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("SmartCard", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$targetPara = $rng.Paragraphs(1)
$insertPoint = $d.Range($targetPara.Range.End, $targetPara.Range.End)

$body1 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">The Other Side</w:t></w:r></w:p><w:p/>'
$insertPoint.InsertXML($pkgHeader + $body1 + $pkgFooter)

# The fragment above carries a trailing empty <w:p/> so that the new
# paragraph gets its own real paragraph mark instead of merging with
# the following ("Resources") paragraph. Track it down and remove it.
$rng2 = $d.Content
$rng2.Find.Execute("The Other Side", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$newPara = $rng2.Paragraphs(1)
$spurious = $newPara.Next()
$spurious.Range.Delete()

# ---------------------------------------------------------------
# Change 2: "Exchanges Listing Guide" + duplicated "Graphics" before
# the existing "Graphics" paragraph; existing paragraph becomes
# "Marketing Materials".
# ---------------------------------------------------------------
$rngG = $d.Content
$rngG.Find.Execute("Graphics", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$graphicsPara = $rngG.Paragraphs(1)
$insertPoint2 = $d.Range($graphicsPara.Range.Start, $graphicsPara.Range.Start)

$body2 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:tab/><w:t xml:space="preserve">Exchanges Listing Guide</w:t></w:r></w:p>' + `
         '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Graphics</w:t></w:r></w:p><w:p/>'
$insertPoint2.InsertXML($pkgHeader + $body2 + $pkgFooter)

# Remove the spurious trailing empty paragraph introduced to terminate
# the duplicated "Graphics" paragraph cleanly.
$rngG2 = $d.Content
$rngG2.Find.Execute("Exchanges Listing Guide", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$elgPara = $rngG2.Paragraphs(1)
$dupGraphicsPara = $elgPara.Next()
$spurious2 = $dupGraphicsPara.Next()
$spurious2.Range.Delete()

# Now turn the ORIGINAL "Graphics" paragraph's run into "Marketing
# Materials" set in Times New Roman (no color override).
$rngOrigGraphics = $dupGraphicsPara.Next()
$runEnd = $rngOrigGraphics.Range.End - 1
$runRange = $d.Range($rngOrigGraphics.Range.Start, $runEnd)
$body3 = '<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">Marketing Materials</w:t></w:r></w:p>'
$runRange.InsertXML($pkgHeader + $body3 + $pkgFooter)

# ---------------------------------------------------------------
# Change 3: "Guide" after "Electrum Wallet"
# ---------------------------------------------------------------
$rngE = $d.Content
$rngE.Find.Execute("Electrum Wallet", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$ewPara = $rngE.Paragraphs(1)
$insertPoint3 = $d.Range($ewPara.Range.End, $ewPara.Range.End)

$body4 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:firstLine="360"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:tab/><w:t>Guide</w:t></w:r></w:p><w:p/>'
$insertPoint3.InsertXML($pkgHeader + $body4 + $pkgFooter)

$rngE2 = $d.Content
$rngE2.Find.Execute("Electrum Wallet", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$ewPara2 = $rngE2.Paragraphs(1)
$guidePara = $ewPara2.Next()
$spurious3 = $guidePara.Next()
$spurious3.Range.Delete()

Write-Host "Edits applied successfully."
